# Weekly update to the Frutilla (strawberry) price sheet for
# "Vega Monumental Concepción": a new price record is inserted as row 267
# (shifting the existing rows 267-273 down to 268-274).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 267; everything below shifts down one row.
$ws.Rows.Item(267).Insert()

# Populate the newly inserted row 267 with the new weekly record.
$ws.Range("A267").Value = 11
$ws.Range("B267").Value = "Vega Monumental Concepción"
$ws.Range("C267").Value = "Bíobío"
$ws.Range("D267").Value = 44595
$ws.Range("E267").Value = 8
$ws.Range("F267").Value = "Fruta"
$ws.Range("G267").Value = 100101
$ws.Range("H267").Value = "Berries"
$ws.Range("I267").Value = 100112025
$ws.Range("J267").Value = "Frutilla"
$ws.Range("K267").Value = "Sin especificar"
$ws.Range("L267").Value = "Primera"
$ws.Range("M267").Value = 250
$ws.Range("N267").Value = 6500
$ws.Range("O267").Value = 7000
$ws.Range("P267").Value = 6740
$ws.Range("Q267").Value = "$/caja 7 kilos"
$ws.Range("R267").Value = "Región del Maule"
$ws.Range("S267").Value = 963
$ws.Range("T267").Value = 7
